$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# RSD sheet: "Max Growth" LPG rate reduced from 5% to 3% (row 33, col C).
# J17 (shared formula "=1+C33") recalculates to 1.03 automatically.
# ---------------------------------------------------------------------------
$rsd = $wb.Worksheets.Item("RSD")
$rsd.Range("C33").Value = 0.03

# ---------------------------------------------------------------------------
# SRV sheet: add a new "LPG" max-growth-rate entry.
# ---------------------------------------------------------------------------
$srv = $wb.Worksheets.Item("SRV")

# Make room for the new summary row directly under the existing ones
# (old row 15 "Max growth rate / Starting value" label row and below
# shift down to make room at row 15).
$srv.Rows.Item(15).Insert()

# Make room for a new raw-data row after the existing Solar row (23)
# and before the footer/header block (old row 24).
$srv.Rows.Item(24).Insert()

# New raw data row (SRV / LPG / 3% max growth / 10% max decline)
$srv.Range("A24").Value = "SRV"
$srv.Range("B24").Value = "LPG"
$srv.Range("C24").Value = 0.03
$srv.Range("D24").Value = 0.1

# New UC summary row referencing the data row above
$srv.Range("B14").Formula = '=_xlfn.TEXTJOIN("_",TRUE,"UC",A24,"MinGrowth",B24)'
$srv.Range("C14").Value = "ACT, GROWTH"
$srv.Range("F14").Value = "SRVLPG"
$srv.Range("G14").Value = "FT*"
$srv.Range("H14").Value = 2021
$srv.Range("I14").Value = "LO"
$srv.Range("J14").Formula = "=1+C24"
$srv.Range("K14").Value = 1
$srv.Range("L14").Formula = "=-D24"
$srv.Range("M14").Value = 5
$srv.Range("N14").Formula = '=_xlfn.TEXTJOIN(" ",TRUE,A24, "maximum growth rate of",B24)'

# ---------------------------------------------------------------------------
# Final UI state: SRV tab active/selected, RSD & SRV selections updated.
# ---------------------------------------------------------------------------
$rsd.Activate()
$rsd.Range("C34").Select()

$srv.Activate()
$srv.Range("K19").Select()
